# Add the 2020 data column (column X) to the 2.1.1 undernourishment sheet,
# mirroring the existing 2019 column (W) for layout/formatting, then fill
# in the new year header and data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number format / font / borders) from the 2019 column (W)
# into the new 2020 column (X) so the new column matches the rest of the table.
$ws.Range("W4:W16").Copy()
$ws.Range("X4:X16").PasteSpecial(-4122)  # xlPasteFormats

# New year header
$ws.Range("X4").Value = 2020

# New 2020 data values (one per country/region row)
$ws.Range("X5").Value = 45.3
$ws.Range("X6").Value = 48.2
$ws.Range("X7").Value = 43.6
$ws.Range("X8").Value = 48.8
$ws.Range("X9").Value = 41.5
$ws.Range("X10").Value = 49.7
$ws.Range("X11").Value = 46.7
$ws.Range("X12").Value = 36.5
$ws.Range("X13").Value = 29.6
$ws.Range("X14").Value = 54.7
$ws.Range("X15").Value = 51.6
$ws.Range("X16").Value = 47.2

# Match the author's final selection in the saved workbook
$ws.Range("AI21").Select() | Out-Null
